# Add a new sale line (item #15) to the DaySale report, pushing the
# totals row and the footer row down by one, and bump the footer
# timestamp to the new save time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert a new row 21, shifting the old total (row 21)
#        and footer (row 22) down to rows 22/23. ------------------------
$ws.Rows.Item(21).Insert()

# Clone row 20 (the last item row) onto the new row 21 so it picks up the
# same cell formatting (fonts/fills/borders/number-formats) and merged
# layout as the rest of the item table, then give it the same row height.
$ws.Range("A20:Q20").Copy($ws.Range("A21:Q21"))
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(21).RowHeight = 25.5

# --- 2. Fill in the new item's data. -------------------------------------
$ws.Range("A21").Value = 15
$ws.Range("C21").Value = "مناديل سولو سحب صغيره"
$ws.Range("N21").Value = "35.00"

# P21 keeps a "0.00" number format but, like the rest of the table, it
# actually holds plain text — flip to a text format while assigning so the
# value isn't coerced into a rounded number, then restore the display
# format used by the other rows.
$pFmt = $ws.Range("P21").NumberFormat()
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "35.0000"
$ws.Range("P21").NumberFormat = $pFmt

# H21 ("0:0"), L21 ("0") and Q21 ("1:0") already carry over correctly from
# row 20's clone, since the new item has the same counters.

# --- 3. Update the totals row (now row 22): add the new line's sell
#        price onto the running total. ------------------------------------
$ws.Range("P22").Value = 298.55

# --- 4. Update the footer row (now row 23) with the new save timestamp. --
$ws.Range("A23").Value = "Sunday, 10 August, 2025 11:02 AM"
